$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values (row => A,B,C,D,E,F) as described by the diff.
$rowData = @{
    3  = @(301, 6, 45, 30, 60, 45)
    4  = @(401, 9, 48, 67, 75, 45)
    5  = @(201, 9, 30, 15, 45, 30)
    6  = @(1201, 2, 10, 10, 10, 10)
    8  = @(701, 3, 90, 45, 97, 15)
    9  = @(1202, 2, 10, 10, 10, 10)
    10 = @(101, 9, 30, 15, 60, 15)
    11 = @(902, 1, 0, 0, 0, 0)
    12 = @(1001, 18, 30, 75, 60, 72)
    15 = @(801, 3, 67, 65, 52, 45)
    16 = @(2, 0, 2, 2, 2, 2)
    17 = @(3, 0, 3, 3, 3, 3)
    18 = @(1101, 0, 15, 30, 30, 0)
    19 = @(1, 0, 2, 2, 2, 2)
    20 = @(502, 0, 4, 0, 0, 0)
    21 = @(802, 0, 4, 5, 4, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
